$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("time costs")

$ws.Range("A5").Value = 41670
$ws.Range("B5").Value = 4.5
$ws.Range("D5").Value = "x"
$ws.Range("E5").Value = "x"
$ws.Range("F5").Value = "x"
$ws.Range("G5").Value = 'Enhancement Data Model & Activity Diagram, Technology check "Appguyver, steroids"'
$ws.Range("H5").Value = "simple custom quizduell app with steroids-technology"

# H5 should pick up the same "center / center" formatting already used by
# D5:F5 (style index 3) rather than column H's default style, so copy the
# format from a cell that already has it instead of minting a new style.
$ws.Range("D2").Copy()
$ws.Range("H5").PasteSpecial(-4122)
$excel.CutCopyMode = $false
